$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "276×6=1656"; New = "859×2=1718" },
    @{ Old = "447×3=1341"; New = "661×7=4627" },
    @{ Old = "721×4=2884"; New = "306×7=2142" },
    @{ Old = "533×7=3731"; New = "489×4=1956" },
    @{ Old = "837×6=5022"; New = "781×3=2343" },
    @{ Old = "155×2=310";  New = "502×2=1004" },
    @{ Old = "354×8=2832"; New = "587×4=2348" },
    @{ Old = "942×6=5652"; New = "359×3=1077" },
    @{ Old = "343×9=3087"; New = "289×4=1156" },
    @{ Old = "960×2=1920"; New = "195×5=975"  },
    @{ Old = "923×7=6461"; New = "441×3=1323" },
    @{ Old = "915×9=8235"; New = "399×9=3591" },
    @{ Old = "540×5=2700"; New = "795×9=7155" },
    @{ Old = "273×5=1365"; New = "991×7=6937" },
    @{ Old = "209×6=1254"; New = "250×7=1750" },
    @{ Old = "301×4=1204"; New = "966×4=3864" },
    @{ Old = "115×4=460";  New = "519×4=2076" },
    @{ Old = "701×8=5608"; New = "163×5=815"  },
    @{ Old = "607×5=3035"; New = "538×8=4304" },
    @{ Old = "494×5=2470"; New = "819×7=5733" },
    @{ Old = "769×4=3076"; New = "230×6=1380" },
    @{ Old = "898×6=5388"; New = "133×2=266"  },
    @{ Old = "921×2=1842"; New = "376×5=1880" },
    @{ Old = "661×9=5949"; New = "493×8=3944" },
    @{ Old = "363×5=1815"; New = "665×7=4655" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2) | Out-Null
}

Write-Output "Replaced $($replacements.Count) cells"
